$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "Plus que 999 heures de jeux"
$ws.Cells.Item(1, 2).Value = "Passez votre première chanson !"
$ws.Cells.Item(2, 1).Value = "Et de deux !"
$ws.Cells.Item(2, 2).Value = "Passez votre première 2 pieds."
$ws.Cells.Item(3, 1).Value = "La force du débutant"
$ws.Cells.Item(3, 2).Value = "Passez votre première 3 pieds."
$ws.Cells.Item(4, 1).Value = "J'ai compris !"
$ws.Cells.Item(4, 2).Value = "Passez votre première 4 pieds."
$ws.Cells.Item(5, 1).Value = "In 8th we trust"
$ws.Cells.Item(5, 2).Value = "Passez votre première 5 pieds."
$ws.Cells.Item(6, 1).Value = "La force tranquille"
$ws.Cells.Item(6, 2).Value = "Passez votre première 6 pieds."
$ws.Cells.Item(7, 1).Value = "Tenace"
$ws.Cells.Item(7, 2).Value = "Passez votre première 7 pieds."
$ws.Cells.Item(8, 1).Value = "Sur le trotoir d'en face"
$ws.Cells.Item(8, 2).Value = "Passez votre première 8 pieds."
$ws.Cells.Item(9, 1).Value = "Bienvenue chez les grands."
$ws.Cells.Item(9, 2).Value = "Passez votre première 9 pieds."
$ws.Cells.Item(10, 1).Value = "Maximum DDR Level"
$ws.Cells.Item(10, 2).Value = "Passez votre première 10 pieds."
$ws.Cells.Item(11, 1).Value = "Un pas de géant"
$ws.Cells.Item(11, 2).Value = "Passez votre première 11 pieds."
$ws.Cells.Item(12, 1).Value = "Premiers saignements"
$ws.Cells.Item(12, 2).Value = "Passez votre première 12 pieds."
$ws.Cells.Item(13, 1).Value = "Maximum ITG Level"
$ws.Cells.Item(13, 2).Value = "Passez votre première 13 pieds."
$ws.Cells.Item(14, 1).Value = "T'es un homme maintenant"
$ws.Cells.Item(14, 2).Value = "Passez votre première 14 pieds."
$ws.Cells.Item(15, 1).Value = "Vers l'infini et l'au-dela !"
$ws.Cells.Item(15, 2).Value = "Passez votre première 15 pieds."
$ws.Cells.Item(16, 1).Value = "Impossible !"
$ws.Cells.Item(16, 2).Value = "Passez votre première 16 pieds."
$ws.Cells.Item(17, 1).Value = "Et maintenant c'est qui le lion ?"
$ws.Cells.Item(17, 2).Value = "Passez votre première 17 pieds."
$ws.Cells.Item(18, 1).Value = "En voilà une longue !"
$ws.Cells.Item(18, 2).Value = "Passez une chanson de plus de 500 steps"
$ws.Cells.Item(19, 1).Value = "Plus de 1000 pas"
$ws.Cells.Item(19, 2).Value = "Passez une chanson de plus de 1000 steps"
$ws.Cells.Item(20, 1).Value = "Interminable"
$ws.Cells.Item(20, 2).Value = "Passez une chanson de plus de 2000 steps"
$ws.Cells.Item(21, 1).Value = "Gloire aux crampes"
$ws.Cells.Item(21, 2).Value = "Passez une chanson de plus de 3000 steps"
$ws.Cells.Item(22, 1).Value = "La fin justifie les moignons"
$ws.Cells.Item(22, 2).Value = "Passez une chanson de plus de 4000 steps"
$ws.Cells.Item(23, 1).Value = "C-C-C-Célimen !"
$ws.Cells.Item(23, 2).Value = "Obtenez votre premier C ou superieur"
$ws.Cells.Item(24, 1).Value = "Bien mais sans plus"
$ws.Cells.Item(24, 2).Value = "Obtenez votre premier B ou superieur"
$ws.Cells.Item(25, 1).Value = "Dans le rythme"
$ws.Cells.Item(25, 2).Value = "Obtenez votre premier A ou superieur"
$ws.Cells.Item(26, 1).Value = "As du rythme"
$ws.Cells.Item(26, 2).Value = "Obtenez votre premier S ou superieur"
$ws.Cells.Item(27, 1).Value = "Sur le podium"
$ws.Cells.Item(27, 2).Value = "Obtenez votre première médaille de bronze"
$ws.Cells.Item(28, 1).Value = "Précision millimétrée"
$ws.Cells.Item(28, 2).Value = "Obtenez votre première médaille d'argent"
$ws.Cells.Item(29, 1).Value = "A deux doigts de la perfection"
$ws.Cells.Item(29, 2).Value = "Obtenez votre première médaille d'or"
$ws.Cells.Item(30, 1).Value = "La perfection"
$ws.Cells.Item(30, 2).Value = "Obtenez votre première médaille de Quad"
$ws.Cells.Item(31, 1).Value = "Aucune tache"
$ws.Cells.Item(31, 2).Value = "Obtenez votre premier full combo"
$ws.Cells.Item(32, 1).Value = "Parcours en or"
$ws.Cells.Item(32, 2).Value = "Obtenez votre premier full excellent combo"
$ws.Cells.Item(33, 1).Value = "Si près de l'or"
$ws.Cells.Item(33, 2).Value = "Obtenez votre premier FC à moins de 10 greats"
$ws.Cells.Item(34, 1).Value = "Score énervant"
$ws.Cells.Item(34, 2).Value = "Obtenez votre premier FEC à moins de 10 ex"
$ws.Cells.Item(35, 1).Value = "Emu aux larmes"
$ws.Cells.Item(35, 2).Value = "Obtenez votre premier FFC"
$ws.Cells.Item(36, 1).Value = "Découverte"
$ws.Cells.Item(36, 2).Value = "Accumulez 500% de score"
$ws.Cells.Item(37, 1).Value = "Interessé"
$ws.Cells.Item(37, 2).Value = "Accumulez plus de 10 000% de score"
$ws.Cells.Item(38, 1).Value = "Geek"
$ws.Cells.Item(38, 2).Value = "Accumulez plus de 50 000% de score"
$ws.Cells.Item(39, 1).Value = "Nolife"
$ws.Cells.Item(39, 2).Value = "Accumulez plus de 100 000% de score"
$ws.Cells.Item(40, 1).Value = "Passioné"
$ws.Cells.Item(40, 2).Value = "Accumulez plus de 200 000% de score"
$ws.Cells.Item(41, 1).Value = "Premiers pas"
$ws.Cells.Item(41, 2).Value = "Accumulez plus de 1000 combos"
$ws.Cells.Item(42, 1).Value = "L'habitude"
$ws.Cells.Item(42, 2).Value = "Accumulez plus de 10 000 combos cumulés"
$ws.Cells.Item(43, 1).Value = "Vétéran du clavier"
$ws.Cells.Item(43, 2).Value = "Accumulez plus de 100 000 combos cumulés"
$ws.Cells.Item(44, 1).Value = "Le million"
$ws.Cells.Item(44, 2).Value = "Atteignez 1 000 000 de combos cumulés"
$ws.Cells.Item(45, 1).Value = "La petite star"
$ws.Cells.Item(45, 2).Value = "Obtenez 5 médailles de bronze ou superieur"
$ws.Cells.Item(46, 1).Value = "Pas de secret"
$ws.Cells.Item(46, 2).Value = "Obtenez 20 médailles de bronze  ou superieur"
$ws.Cells.Item(47, 1).Value = "Fini de rire"
$ws.Cells.Item(47, 2).Value = "Obtenez 100 médailles de bronze  ou superieur"
$ws.Cells.Item(48, 1).Value = "La starlette"
$ws.Cells.Item(48, 2).Value = "Obtenez 5 médailles d'argent  ou superieur"
$ws.Cells.Item(49, 1).Value = "La maitrise"
$ws.Cells.Item(49, 2).Value = "Obtenez 20 médailles d'argent  ou superieur"
$ws.Cells.Item(50, 1).Value = "VIP de Cublast"
$ws.Cells.Item(50, 2).Value = "Obtenez 100 médailles d'argent  ou superieur"
$ws.Cells.Item(51, 1).Value = "La superstar"
$ws.Cells.Item(51, 2).Value = "Obtenez 5 médaille d'or  ou superieur"
$ws.Cells.Item(52, 1).Value = "La consécration"
$ws.Cells.Item(52, 2).Value = "Obtenez 10 médailles d'or  ou superieur"
$ws.Cells.Item(53, 1).Value = "Plus de doute"
$ws.Cells.Item(53, 2).Value = "Obtenez 50 médailles d'or  ou superieur"
$ws.Cells.Item(54, 1).Value = "La légende"
$ws.Cells.Item(54, 2).Value = "Obtenez 2 médaille de Quad"
$ws.Cells.Item(55, 1).Value = "Continuer la légende"
$ws.Cells.Item(55, 2).Value = "Obtenez 3 médailles de Quad"
$ws.Cells.Item(56, 1).Value = "Ecrire la légende"
$ws.Cells.Item(56, 2).Value = "Obtenez 5 médailles de Quad"
$ws.Cells.Item(57, 1).Value = "Mélomane"
$ws.Cells.Item(57, 2).Value = "Clear 50 chansons différentes"
$ws.Cells.Item(58, 1).Value = "Curieux"
$ws.Cells.Item(58, 2).Value = "Clear 250 chansons différentes"
$ws.Cells.Item(59, 1).Value = "4Go sur le disque dur"
$ws.Cells.Item(59, 2).Value = "Clear 600 chanson différentes"
$ws.Cells.Item(60, 2).Value = "Battez le boss 1"
$ws.Cells.Item(61, 2).Value = "Battez le boss 2"
$ws.Cells.Item(62, 2).Value = "Battez le boss 3"
$ws.Cells.Item(63, 2).Value = "Battez le boss 4"
$ws.Cells.Item(64, 2).Value = "Battez le boss 5"
$ws.Cells.Item(65, 2).Value = "Battez le boss 6"
$ws.Cells.Item(66, 2).Value = "Battez le boss 7"
$ws.Cells.Item(67, 2).Value = "Battez le boss 8"
$ws.Cells.Item(68, 1).Value = "Les bonnes choses ont une fin"
$ws.Cells.Item(68, 2).Value = "Affrontez le boss de fin et finissez Cublast"
$ws.Cells.Item(69, 1).Value = "Les mauvaises choses ont une fin"
$ws.Cells.Item(69, 2).Value = "Retourner affronter le boss de fin pour la 2eme fois"
$ws.Cells.Item(70, 1).Value = "Toutes les choses ont une fin"
$ws.Cells.Item(70, 2).Value = "Obtenez la medaille de Quad du boss de fin"
$ws.Cells.Item(71, 1).Value = "Sherlock"
$ws.Cells.Item(71, 2).Value = "Debloquez toutes les chansons cachées"
$ws.Cells.Item(72, 1).Value = "Challenger"
$ws.Cells.Item(72, 2).Value = "Finissez toutes les chansons challenges"
$ws.Cells.Item(73, 1).Value = "La sauvegarde du forum "
$ws.Cells.Item(73, 2).Value = "Obtenez toutes les médailles de bronze dans le mode Story"
$ws.Cells.Item(74, 1).Value = "La paix du forum"
$ws.Cells.Item(74, 2).Value = "Obtenez toutes les médailles d'argent dans le mode Story"
$ws.Cells.Item(75, 1).Value = "Le règne sur le forum"
$ws.Cells.Item(75, 2).Value = "Obtenez toutes les médailles d'or dans le mode Story"
$ws.Cells.Item(76, 1).Value = "Le maitre du forum"
$ws.Cells.Item(76, 2).Value = "Obtenez toutes les médailles de Quad dans le mode Story"
$ws.Cells.Item(77, 1).Value = "Joueur social"
$ws.Cells.Item(77, 2).Value = "Affrontez 1 ami en mode online en même temps"
$ws.Cells.Item(78, 1).Value = "Ramenez vous les potes"
$ws.Cells.Item(78, 2).Value = "Affrontez 3 amis en mode online en même temps"
$ws.Cells.Item(79, 1).Value = "Faire passer le fun"
$ws.Cells.Item(79, 2).Value = "Récuperez 5 profiles de joueurs différents"
$ws.Cells.Item(80, 1).Value = "La communauté"
$ws.Cells.Item(80, 2).Value = "Récuperez 10 profiles de joueurs différents"
$ws.Cells.Item(81, 1).Value = "Mise à jour"
$ws.Cells.Item(81, 2).Value = "Mettez à jour le profil d'un joueur que vous aviez déjà"
$ws.Cells.Item(82, 1).Value = "Le gout du combat"
$ws.Cells.Item(82, 2).Value = "Battez le score d'un de vos amis en mode solo"
$ws.Cells.Item(83, 1).Value = "Le gout de la compétition"
$ws.Cells.Item(83, 2).Value = "Battez 10 scores de vos amis en mode solo"
$ws.Cells.Item(84, 1).Value = "Toutes les avoir ne suffit pas"
$ws.Cells.Item(84, 2).Value = "Gagnez un match contre un ami alors qu'il a fait un full combo"

$ws.Range("B73").Select()
